$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: label + same style as the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data cells F2:F105: per-row "time_taken" timestamps recorded for each panel row
$timeTaken = @("2021-10-05 10:52:28.913140","2021-10-05 10:52:28.913151","2021-10-05 10:52:28.913155","2021-10-05 10:52:28.913157","2021-10-05 10:52:28.913160","2021-10-05 10:52:28.913163","2021-10-05 10:52:28.913165","2021-10-05 10:52:28.913168","2021-10-05 10:52:28.913170","2021-10-05 10:52:28.913173","2021-10-05 10:52:28.913176","2021-10-05 10:52:28.913178","2021-10-05 10:52:28.913180","2021-10-05 10:52:28.913183","2021-10-05 10:52:28.913186","2021-10-05 10:52:28.913188","2021-10-05 10:52:28.913191","2021-10-05 10:52:28.913193","2021-10-05 10:52:28.913196","2021-10-05 10:52:28.913198","2021-10-05 10:52:28.913201","2021-10-05 10:52:28.913203","2021-10-05 10:52:28.913206","2021-10-05 10:52:28.913209","2021-10-05 10:52:28.913212","2021-10-05 10:52:28.913215","2021-10-05 10:52:28.913218","2021-10-05 10:52:28.913220","2021-10-05 10:52:28.913223","2021-10-05 10:52:28.913226","2021-10-05 10:52:28.913229","2021-10-05 10:52:28.913231","2021-10-05 10:52:28.913235","2021-10-05 10:52:28.913237","2021-10-05 10:52:28.913240","2021-10-05 10:52:28.913242","2021-10-05 10:52:28.913245","2021-10-05 10:52:28.913247","2021-10-05 10:52:28.913250","2021-10-05 10:52:28.913253","2021-10-05 10:52:28.913255","2021-10-05 10:52:28.913258","2021-10-05 10:52:28.913260","2021-10-05 10:52:28.913263","2021-10-05 10:52:28.913265","2021-10-05 10:52:28.913268","2021-10-05 10:52:28.913271","2021-10-05 10:52:28.913273","2021-10-05 10:52:28.913276","2021-10-05 10:52:28.913278","2021-10-05 10:52:28.913281","2021-10-05 10:52:28.913283","2021-10-05 10:52:28.913286","2021-10-05 10:52:28.913289","2021-10-05 10:52:28.913292","2021-10-05 10:52:28.913294","2021-10-05 10:52:28.913296","2021-10-05 10:52:28.913299","2021-10-05 10:52:28.913302","2021-10-05 10:52:28.913304","2021-10-05 10:52:28.913307","2021-10-05 10:52:28.913309","2021-10-05 10:52:28.913312","2021-10-05 10:52:28.913315","2021-10-05 10:52:28.913319","2021-10-05 10:52:28.913321","2021-10-05 10:52:28.913324","2021-10-05 10:52:28.913327","2021-10-05 10:52:28.913330","2021-10-05 10:52:28.913332","2021-10-05 10:52:28.913335","2021-10-05 10:52:28.913337","2021-10-05 10:52:28.913340","2021-10-05 10:52:28.913343","2021-10-05 10:52:28.913346","2021-10-05 10:52:28.913349","2021-10-05 10:52:28.913354","2021-10-05 10:52:28.913357","2021-10-05 10:52:28.913360","2021-10-05 10:52:28.913362","2021-10-05 10:52:28.913365","2021-10-05 10:52:28.913367","2021-10-05 10:52:28.913370","2021-10-05 10:52:28.913372","2021-10-05 10:52:28.913375","2021-10-05 10:52:28.913378","2021-10-05 10:52:28.913380","2021-10-05 10:52:28.913383","2021-10-05 10:52:28.913386","2021-10-05 10:52:28.913388","2021-10-05 10:52:28.913391","2021-10-05 10:52:28.913393","2021-10-05 10:52:28.913397","2021-10-05 10:52:28.913400","2021-10-05 10:52:28.913403","2021-10-05 10:52:28.913406","2021-10-05 10:52:28.913409","2021-10-05 10:52:28.913411","2021-10-05 10:52:28.913414","2021-10-05 10:52:28.913417","2021-10-05 10:52:28.913419","2021-10-05 10:52:28.913422","2021-10-05 10:52:28.913425","2021-10-05 10:52:28.913428")

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "Applied time_taken column to F1:F105"
